{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) The \"git init\" list item is currently split across three runs\n//    (with spell-check proofErr markers wrapping \"init\"). Re-set the\n//    paragraph's own text in place so Word collapses it back down to a\n//    single run (dropping the now-stale proofErr markers) while leaving\n//    the visible text \u2013 and the run's language formatting \u2013 unchanged.\nconst initPara = items.find((p) => p.text.indexOf(\"git init\") === 0);\nif (initPara) {\n  initPara.getRange().insertText(initPara.text, \"Replace\");\n}\n\n// 2) The document ends with an empty list-style paragraph. Fill it in\n//    with the new \"git branch -d\" command description, matching the\n//    run-level formatting (en-US language) used by its sibling bullets.\nconst lastPara = items[items.length - 1];\nconst newCommandOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>git branch -d (branch name) \\u2013 This command will delete the particular branch from the local repository.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nlastPara.insertOoxml(newCommandOoxml, \"End\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The \"git init\" bullet is split across three runs (with stale\n#    spell-check proofErr markers wrapping \"init\"). Running Find/Replace\n#    over the paragraph's own full text re-types it as a single run,\n#    dropping the proofErr markers while keeping the run's formatting\n#    (en-US language) and the paragraph's own identity untouched.\n$enDash = [char]8211\n$initText = \"git init \" + $enDash + \" Creates a new Git repository in the path where git command is initialized.\"\n$findRange = $d.Content\n$findRange.Find.Execute($initText, $true, $false, $false, $false, $false, $true, 1, $false, $initText, 2)\n\n# 2) The document ends with an empty list-style paragraph. Copy the\n#    formatting of the previous bullet (so the new run picks up the same\n#    en-US language run formatting) and then replace its text with the\n#    new \"git branch -d\" command description.\n$count = $d.Paragraphs.Count\n$srcPara = $d.Paragraphs.Item($count - 1)\n$lastPara = $d.Paragraphs.Item($count)\n\n$lastPara.Range.FormattedText = $srcPara.Range.FormattedText\n$newCommandText = \"git branch -d (branch name) \" + $enDash + \" This command will delete the particular branch from the local repository.\"\n$lastPara.Range.Text = $newCommandText\n"}
